# Doing Updates for Financials
# Update the yearly financial figures (Income Statement, Balance Sheet,
# Cash Flow Statement) on the ELLO sheet to the restated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income Statement - Total Revenue (row 8)
$ws.Range("D8").Value = 15300
$ws.Range("E8").Value = 13100
$ws.Range("F8").Value = 14000
$ws.Range("G8").Value = 17700
$ws.Range("H8").Value = 14600
$ws.Range("I8").Value = 10000
$ws.Range("J8").Value = 6900

$ws.Range("D9").Value = 7900
$ws.Range("E9").Value = 7300
$ws.Range("F9").Value = 7900
$ws.Range("G9").Value = 9600
$ws.Range("H9").Value = 7200
$ws.Range("I9").Value = 5200
$ws.Range("J9").Value = 7100

# Gross Profit (row 10)
$ws.Range("D10").Value = 7400
$ws.Range("E10").Value = 5800
$ws.Range("F10").Value = 6100
$ws.Range("G10").Value = 8100
$ws.Range("H10").Value = 7400
$ws.Range("I10").Value = 4700
$ws.Range("J10").Value = -200

$ws.Range("D12").Value = 3100
$ws.Range("E12").Value = 2500

$ws.Range("E14").Value = 2000
$ws.Range("G14").Value = -4500
$ws.Range("H14").Value = -11500
$ws.Range("I14").Value = 800

$ws.Range("D17").Value = 12000
$ws.Range("E17").Value = 10400
$ws.Range("F17").Value = 9100
$ws.Range("G17").Value = 6200
$ws.Range("I17").Value = 9300
$ws.Range("J17").Value = 7000

$ws.Range("D18").Value = 3300
$ws.Range("E18").Value = 2700
$ws.Range("F18").Value = 4800
$ws.Range("G18").Value = 11500
$ws.Range("H18").Value = 14300

$ws.Range("D20").Value = -9500
$ws.Range("E20").Value = -2200
$ws.Range("F20").Value = 2700
$ws.Range("G20").Value = -2900
$ws.Range("H20").Value = -1300
$ws.Range("I20").Value = -3100
$ws.Range("J20").Value = -1500

$ws.Range("E21").Value = 5500
$ws.Range("F21").Value = 12500
$ws.Range("G21").Value = 14700
$ws.Range("H21").Value = 17500

$ws.Range("H22").Value = 1400

$ws.Range("D23").Value = -7000
$ws.Range("F23").Value = 7100
$ws.Range("G23").Value = 7700
$ws.Range("H23").Value = 11600
$ws.Range("I23").Value = -3500
$ws.Range("J23").Value = -2200

$ws.Range("E24").Value = 600
$ws.Range("I24").Value = -1100
$ws.Range("J24").Value = -1100

$ws.Range("D26").Value = -7500
$ws.Range("F26").Value = 9100
$ws.Range("G26").Value = 7500
$ws.Range("H26").Value = 11300
$ws.Range("I26").Value = -2400

$ws.Range("D27").Value = -6900
$ws.Range("F27").Value = 9400
$ws.Range("G27").Value = 7500
$ws.Range("H27").Value = 11300
$ws.Range("I27").Value = -2400

$ws.Range("D32").Value = 9500
$ws.Range("E32").Value = 2200
$ws.Range("F32").Value = -2700
$ws.Range("G32").Value = 2900
$ws.Range("H32").Value = 1300
$ws.Range("I32").Value = 3100
$ws.Range("J32").Value = 1500

$ws.Range("D33").Value = -6900
$ws.Range("F33").Value = 9400
$ws.Range("G33").Value = 7500
$ws.Range("H33").Value = 11300
$ws.Range("I33").Value = -2400

$ws.Range("D35").Value = -6900
$ws.Range("F35").Value = 9400
$ws.Range("G35").Value = 7500
$ws.Range("H35").Value = 11300
$ws.Range("I35").Value = -2400

$ws.Range("D41").Value = 26900
$ws.Range("E41").Value = 25200
$ws.Range("F41").Value = 21000
$ws.Range("G41").Value = 17700
$ws.Range("H41").Value = 8100
$ws.Range("I41").Value = 37400
$ws.Range("J41").Value = 32400

$ws.Range("D42").Value = 3800
$ws.Range("F42").Value = 7300
$ws.Range("G42").Value = 8600

$ws.Range("D43").Value = 12800
$ws.Range("E43").Value = 9800
$ws.Range("F43").Value = 8500
$ws.Range("G43").Value = 6300
$ws.Range("H43").Value = 4100
$ws.Range("I43").Value = 4200
$ws.Range("J43").Value = 14600

$ws.Range("D45").Value = 4300
$ws.Range("E45").Value = 800
$ws.Range("G45").Value = 900
$ws.Range("H45").Value = 13000
$ws.Range("I45").Value = 15900
$ws.Range("J45").Value = 58100

$ws.Range("D46").Value = 47800
$ws.Range("E46").Value = 37000
$ws.Range("F46").Value = 37600
$ws.Range("G46").Value = 33500
$ws.Range("H46").Value = 25300
$ws.Range("I46").Value = 57400
$ws.Range("J46").Value = 69100

$ws.Range("D47").Value = 73800
$ws.Range("E47").Value = 38800
$ws.Range("F47").Value = 44500
$ws.Range("G47").Value = 34400
$ws.Range("H47").Value = 29600
$ws.Range("I47").Value = 22900
$ws.Range("J47").Value = 14800

$ws.Range("D48").Value = 88500
$ws.Range("E48").Value = 82200
$ws.Range("F48").Value = 88600
$ws.Range("G48").Value = 104900
$ws.Range("H48").Value = 105100
$ws.Range("I48").Value = 60400
$ws.Range("J48").Value = 52600

$ws.Range("D49").Value = 6200

$ws.Range("D52").Value = 6100
$ws.Range("E52").Value = 8500
$ws.Range("F52").Value = 9200
$ws.Range("G52").Value = 5800
$ws.Range("H52").Value = 4800
$ws.Range("I52").Value = 3600
$ws.Range("J52").Value = 5900

$ws.Range("D54").Value = 222300
$ws.Range("E54").Value = 166600
$ws.Range("F54").Value = 179900
$ws.Range("G54").Value = 178500
$ws.Range("H54").Value = 164900
$ws.Range("I54").Value = 144400
$ws.Range("J54").Value = 141800

$ws.Range("D57").Value = 1500
$ws.Range("E57").Value = 1800
$ws.Range("H57").Value = 3000
$ws.Range("I57").Value = 12700
$ws.Range("J57").Value = 28100

$ws.Range("D58").Value = 8700
$ws.Range("E58").Value = 6600
$ws.Range("F58").Value = 6700
$ws.Range("G58").Value = 6200
$ws.Range("H58").Value = 21800
$ws.Range("I58").Value = 7900
$ws.Range("J58").Value = 13600

$ws.Range("D59").Value = 2500
$ws.Range("E59").Value = 3500
$ws.Range("F59").Value = 3600
$ws.Range("G59").Value = 4600
$ws.Range("H59").Value = 5400
$ws.Range("I59").Value = 5400
$ws.Range("J59").Value = 8000

$ws.Range("D60").Value = 12700
$ws.Range("E60").Value = 11800
$ws.Range("F60").Value = 11300
$ws.Range("G60").Value = 12300
$ws.Range("H60").Value = 30200
$ws.Range("I60").Value = 26100
$ws.Range("J60").Value = 33300

$ws.Range("D61").Value = 110800
$ws.Range("E61").Value = 56100
$ws.Range("F61").Value = 59300
$ws.Range("G61").Value = 55800
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 20800
$ws.Range("J61").Value = 12600

$ws.Range("D62").Value = 11800
$ws.Range("E62").Value = 3900
$ws.Range("F62").Value = 3700
$ws.Range("G62").Value = 4800
$ws.Range("H62").Value = 2700
$ws.Range("I62").Value = 4300
$ws.Range("J62").Value = 1600

$ws.Range("D66").Value = 134000
$ws.Range("E66").Value = 71100
$ws.Range("F66").Value = 74000
$ws.Range("G66").Value = 72900
$ws.Range("H66").Value = 53000
$ws.Range("I66").Value = 51200
$ws.Range("J66").Value = 47500

$ws.Range("D72").Value = 2300
$ws.Range("E72").Value = 9500
$ws.Range("F72").Value = -9000
$ws.Range("G72").Value = -9500
$ws.Range("H72").Value = -3200
$ws.Range("I72").Value = -21300
$ws.Range("J72").Value = -20700

$ws.Range("D76").Value = 88200
$ws.Range("E76").Value = 95400
$ws.Range("F76").Value = 105800
$ws.Range("G76").Value = 105600
$ws.Range("H76").Value = 111900
$ws.Range("I76").Value = 93200
$ws.Range("J76").Value = 94300

$ws.Range("D81").Value = -6900
$ws.Range("F81").Value = 9400
$ws.Range("G81").Value = 7500
$ws.Range("H81").Value = 11300
$ws.Range("I81").Value = -2400

$ws.Range("D83").Value = 5100
$ws.Range("E83").Value = 4900
$ws.Range("F83").Value = 5000
$ws.Range("G83").Value = 6100
$ws.Range("H83").Value = 4500
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2000

$ws.Range("D89").Value = 2600
$ws.Range("E89").Value = 8200
$ws.Range("F89").Value = 5200
$ws.Range("G89").Value = 3700
$ws.Range("H89").Value = 7200
$ws.Range("I89").Value = 6600
$ws.Range("J89").Value = -3400

$ws.Range("D91").Value = -8500
$ws.Range("E91").Value = -5700
$ws.Range("H91").Value = -10300
$ws.Range("J91").Value = -28000

$ws.Range("D94").Value = -30700
$ws.Range("E94").Value = 600
$ws.Range("F94").Value = -4600
$ws.Range("G94").Value = -18000
$ws.Range("H94").Value = -48000
$ws.Range("I94").Value = -2100
$ws.Range("J94").Value = -70500

$ws.Range("E96").Value = -2400

$ws.Range("D100").Value = 33300
$ws.Range("E100").Value = -2700
$ws.Range("F100").Value = 4300
$ws.Range("G100").Value = 28000
$ws.Range("H100").Value = 11100
$ws.Range("J100").Value = 23600

$ws.Range("D101").Value = -3500
$ws.Range("G101").Value = -4100
$ws.Range("J101").Value = -3200

$ws.Range("E102").Value = 5900
$ws.Range("F102").Value = 4700
$ws.Range("G102").Value = 9600
$ws.Range("H102").Value = -29200
$ws.Range("I102").Value = 4900
$ws.Range("J102").Value = -53500

